# Update code and metadata
# - Reorder the variable/description rows (lat/lon moved up, study_ID moved
#   down) to match the data dictionary produced by the updated analysis code.
# - Document new/renamed variables used by the updated model: PC.predation
#   (with a mixed-formatting description), density, and the NPP-related
#   helper columns (metric, NPP.proxy, NPP.scale, NPP.scale2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing rows, re-ordered / re-filled in their new positions ---
$ws.Cells.Item(1,1).Value = "Variable"
$ws.Cells.Item(1,2).Value = "Description"
$ws.Cells.Item(2,1).Value = "FW_name"
$ws.Cells.Item(2,2).Value = "A unique identifier for each individual food web"
$ws.Cells.Item(3,1).Value = "temperature_C"
$ws.Cells.Item(3,2).Value = "Environmental temperature (°C)"
$ws.Cells.Item(4,1).Value = "lat"
$ws.Cells.Item(4,2).Value = "Latitude coordinate for food web location"
$ws.Cells.Item(5,1).Value = "lon"
$ws.Cells.Item(5,2).Value = "Longitude coordinate for food web location"
$ws.Cells.Item(6,1).Value = "ecosystem.type"
$ws.Cells.Item(6,2).Value = "Which type of ecosystem the data come from"
$ws.Cells.Item(7,1).Value = "sampling.start.year"
$ws.Cells.Item(7,2).Value = "Year in which first sampling took place"
$ws.Cells.Item(8,1).Value = "sampling.end.year"
$ws.Cells.Item(8,2).Value = "Year in which final sampling took place"
$ws.Cells.Item(9,1).Value = "study_ID"
$ws.Cells.Item(9,2).Value = "A unique identifier for each study from which data are derived"
$ws.Cells.Item(10,1).Value = "stability"
$ws.Cells.Item(10,2).Value = "The leading eigenvalue of the Jacobian matrix"
$ws.Cells.Item(11,1).Value = "flux"
$ws.Cells.Item(11,2).Value = "The total summed energy flux through the entire food web"
$ws.Cells.Item(12,1).Value = "second.consumption"
$ws.Cells.Item(12,2).Value = "Total energy outflux from all consumer nodes"
$ws.Cells.Item(13,1).Value = "prim.consumption"
$ws.Cells.Item(13,2).Value = "Total energy outflux from all basal resource nodes"
$ws.Cells.Item(15,1).Value = "S"
$ws.Cells.Item(15,2).Value = "Taxa richness of the food web"
$ws.Cells.Item(17,1).Value = "MaxTL"
$ws.Cells.Item(17,2).Value = "Maximum trophic level"
$ws.Cells.Item(18,1).Value = "sim.sec.cons"
$ws.Cells.Item(18,2).Value = "Mean trophic similarity of secondary conumers "
$ws.Cells.Item(19,1).Value = "sim.prim.cons"
$ws.Cells.Item(19,2).Value = "Mean trophic similarity of primary conumers "
$ws.Cells.Item(20,1).Value = "sim.total"
$ws.Cells.Item(20,2).Value = "Mean trophic similarity of the entire food web"

# --- New variable labels (column A) for the newly documented rows ---
$ws.Cells.Item(14,1).Value = "PC.predation"
$ws.Cells.Item(16,1).Value = "density"
$ws.Cells.Item(21,1).Value = "metric"
$ws.Cells.Item(22,1).Value = "NPP.proxy"
$ws.Cells.Item(23,1).Value = "NPP.scale"
$ws.Cells.Item(24,1).Value = "NPP.scale2"

# --- New variable descriptions (column B) ---
$ws.Cells.Item(21,2).Value = "measure of NPP (NDVI or chlorophyll-a)"
$ws.Cells.Item(22,2).Value = "Raw NPP value derived from metric"
$ws.Cells.Item(16,2).Value = "Estimated consumer density (abundance/area) of the food web"
$ws.Cells.Item(23,2).Value = "logit-transformed NPP.proxy"
$ws.Cells.Item(24,2).Value = "NPP.scale^2 (quadratic term)"

# PC.predation description is added last, with "Per capita" in italics and
# the remainder of the sentence in normal (upright) text.
$ws.Cells.Item(14,2).Value = "Per capita predation rate (prey outflux / prey biomass)"
$ws.Cells.Item(14,2).Font.Italic = $true
$ws.Cells.Item(14,2).Characters(11, 46).Font.Italic = $false

# Column A best-fits to the new (slightly longer) set of labels.
$ws.Columns("A:A").AutoFit()

# Leave the selection on the last-edited cell, as in the source workbook.
$ws.Cells.Item(18,2).Select()
